$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 35; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = "B-" + $cell.Value2
}
